$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text type (avoids Excel auto-numeric
# conversion for numeric-looking strings), without leaving any permanent
# style/number-format change behind.
function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Cells.Item(2, 4).Value = '27.970.66'
$ws.Cells.Item(2, 5).Value = '  +1.18%  '

$ws.Cells.Item(3, 4).Value = '1.642.46'
$ws.Cells.Item(3, 5).Value = '  +0.53%  '

$ws.Cells.Item(4, 5).Value = '  +0.04%  '

Set-TextValue $ws.Cells.Item(6, 4) '0.525'
$ws.Cells.Item(6, 5).Value = '  +0.45%  '

$ws.Cells.Item(7, 5).Value = '  +0.01%  '

Set-TextValue $ws.Cells.Item(8, 4) '23.56'
$ws.Cells.Item(8, 5).Value = '  +1.58%  '

Set-TextValue $ws.Cells.Item(9, 4) '0.259'
$ws.Cells.Item(9, 5).Value = '  -1.67%  '

$ws.Cells.Item(10, 5).Value = '  +0.61%  '

Set-TextValue $ws.Cells.Item(11, 4) '0.0884'
$ws.Cells.Item(11, 5).Value = '  +2.49%  '

$ws.Cells.Item(12, 4).Value = '1.874.93'
$ws.Cells.Item(12, 5).Value = '  +0.48%  '

$ws.Cells.Item(13, 4).Value = '1.641.42'
$ws.Cells.Item(13, 5).Value = '  +0.47%  '

$ws.Cells.Item(14, 5).Value = '  +1.13%  '

Set-TextValue $ws.Cells.Item(15, 4) '0.573'
$ws.Cells.Item(15, 5).Value = '  +2.15%  '

Set-TextValue $ws.Cells.Item(16, 4) '65.54'
$ws.Cells.Item(16, 5).Value = '  +0.56%  '

$ws.Cells.Item(17, 4).Value = '27.971.26'
$ws.Cells.Item(17, 5).Value = '  +1.22%  '

Set-TextValue $ws.Cells.Item(18, 4) '233.16'
$ws.Cells.Item(18, 5).Value = '  +1.42%  '

$ws.Cells.Item(19, 5).Value = '  +0.57%  '

Set-TextValue $ws.Cells.Item(20, 4) '7.61'
$ws.Cells.Item(20, 5).Value = '  +0.05%  '

$ws.Cells.Item(21, 5).Value = '  +0.12%  '

Set-TextValue $ws.Cells.Item(22, 4) '10.58'
$ws.Cells.Item(22, 5).Value = '  -0.81%  '

Set-TextValue $ws.Cells.Item(23, 4) '4.38'
$ws.Cells.Item(23, 5).Value = '  -0.09%  '

$ws.Cells.Item(24, 5).Value = '  -2.88%  '

Set-TextValue $ws.Cells.Item(25, 4) '152.91'
$ws.Cells.Item(25, 5).Value = '  +2.22%  '

$ws.Cells.Item(26, 5).Value = '  +0.44%  '

$ws.Cells.Item(27, 5).Value = '  +0.39%  '

$ws.Cells.Item(28, 5).Value = '  +0.24%  '

$ws.Cells.Item(29, 5).Value = '  +0.05%  '

$ws.Cells.Item(30, 5).Value = '  +0.58%  '

$ws.Cells.Item(31, 5).Value = '  +1.00%  '

$ws.Cells.Item(32, 5).Value = '  +3.74%  '

$ws.Cells.Item(33, 5).Value = '  +0.39%  '

$ws.Cells.Item(34, 4).Value = '1.410.14'
$ws.Cells.Item(34, 5).Value = '  -3.87%  '

Set-TextValue $ws.Cells.Item(35, 4) '1.59'
$ws.Cells.Item(35, 5).Value = '  +2.23%  '

$ws.Cells.Item(36, 5).Value = '  +1.83%  '

$ws.Cells.Item(37, 2).Value = 'ImmutableX'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Cells.Item(37, 4) '0.566'
$ws.Cells.Item(37, 5).Value = '  +1.62%  '

$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Cells.Item(38, 4) '0.0169'
$ws.Cells.Item(38, 5).Value = '  +1.30%  '

Set-TextValue $ws.Cells.Item(39, 4) '0.881'
$ws.Cells.Item(39, 5).Value = '  +0.35%  '

$ws.Cells.Item(40, 5).Value = '  +0.34%  '

$ws.Cells.Item(41, 5).Value = '  +1.35%  '

$ws.Cells.Item(42, 5).Value = '  +0.08%  '

$ws.Cells.Item(43, 5).Value = '  +7.26%  '

Set-TextValue $ws.Cells.Item(44, 4) '67.22'
$ws.Cells.Item(44, 5).Value = '  -2.45%  '

$ws.Cells.Item(45, 5).Value = '  +3.37%  '

$ws.Cells.Item(46, 5).Value = '  +0.18%  '

$ws.Cells.Item(47, 4).Value = '1.783.92'
$ws.Cells.Item(47, 5).Value = '  +0.44%  '

Set-TextValue $ws.Cells.Item(48, 4) '88.11'
$ws.Cells.Item(48, 5).Value = '  +0.50%  '

$ws.Cells.Item(49, 5).Value = '  +0.50%  '

$ws.Cells.Item(50, 5).Value = '  +0.41%  '

Set-TextValue $ws.Cells.Item(51, 4) '7.60'
$ws.Cells.Item(51, 5).Value = '  -0.94%  '
